$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 10
$ws.Range("F4").Value = 402
$ws.Range("F5").Value = 202
$ws.Range("F6").Value = 809
$ws.Range("F7").Value = 112
$ws.Range("F8").Value = 10309
$ws.Range("F9").Value = 57
$ws.Range("F10").Value = 3555
$ws.Range("F12").Value = 2461
$ws.Range("F13").Value = 40
$ws.Range("F14").Value = 2839
$ws.Range("F16").Value = 510
$ws.Range("F17").Value = 2190
$ws.Range("F18").Value = 46
$ws.Range("F20").Value = 29
$ws.Range("F21").Value = 394
$ws.Range("F22").Value = 22
$ws.Range("F23").Value = 156
$ws.Range("F24").Value = 319
$ws.Range("F25").Value = 276
$ws.Range("F26").Value = 237
$ws.Range("F28").Value = 1328
$ws.Range("F29").Value = 16
$ws.Range("F31").Value = 106
$ws.Range("F32").Value = 132
$ws.Range("F34").Value = 3858
$ws.Range("F35").Value = 3257
$ws.Range("F40").Value = 6
$ws.Range("F42").Value = 110
$ws.Range("F43").Value = 111
$ws.Range("F44").Value = 74
$ws.Range("F45").Value = 29
$ws.Range("F46").Value = 44
$ws.Range("F47").Value = 14

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 10
$ws.Range("F4").Value = 181
$ws.Range("F15").Value = 37

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 755
$ws.Range("F3").Value = 993
$ws.Range("F4").Value = 130
$ws.Range("F5").Value = 2072

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 755
$ws.Range("F3").Value = 993
$ws.Range("F4").Value = 130
$ws.Range("F5").Value = 10
$ws.Range("F6").Value = 402
$ws.Range("F8").Value = 202
$ws.Range("F9").Value = 809
$ws.Range("F10").Value = 112
$ws.Range("F11").Value = 10309
$ws.Range("F12").Value = 57
$ws.Range("F13").Value = 3555
$ws.Range("F15").Value = 2461
$ws.Range("F16").Value = 40
$ws.Range("F18").Value = 510
$ws.Range("F19").Value = 2190
$ws.Range("F20").Value = 46
$ws.Range("F22").Value = 29
$ws.Range("F23").Value = 394
$ws.Range("F24").Value = 156
$ws.Range("F25").Value = 319
$ws.Range("F26").Value = 237
$ws.Range("F27").Value = 1328
$ws.Range("F28").Value = 16
$ws.Range("F30").Value = 106
$ws.Range("F31").Value = 132
$ws.Range("F36").Value = 3258
$ws.Range("F41").Value = 6
$ws.Range("F43").Value = 37
$ws.Range("F45").Value = 110
$ws.Range("F46").Value = 74
$ws.Range("F47").Value = 29
$ws.Range("F48").Value = 14

